$d = $word.ActiveDocument

# Locate the bulleted list item "Oldest and the simplest data storage."
# (the last bullet of the first BLOCK STORAGE section) -- the new bullet
# needs to be inserted directly after it, before the following blank
# paragraph / VOLUME heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Oldest and the simplest data storage.") {
        $target = $p
        break
    }
}

# Put the insertion point at the end of that paragraph (after the trailing
# period, before the paragraph mark) and create a brand new paragraph right
# after it.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# The newly created paragraph is now the one right after $target.
$inserted = $target.Next()
$inserted.Range.Text = "Would not contain meta data"

# Give it the same list formatting (ListParagraph style, numId 1 / ilvl 0)
# as the rest of the bulleted list it belongs to.
$inserted.Style = $target.Style
$inserted.Range.ListFormat.ListTemplate = $target.Range.ListFormat.ListTemplate
$inserted.Range.ListFormat.ListLevelNumber = $target.Range.ListFormat.ListLevelNumber
